$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://instrument-orugie.ru/catalog/Sadovaya-tehnika/Motobloki/"

# Register the new hyperlink (this also writes the URL text into B2, as a new
# shared string, and wires up the external relationship for it).
$ws.Hyperlinks.Add($ws.Range("B2"), $url)

# Hyperlinks.Add stamps its own generic "Hyperlink" cell style; re-apply the
# same look already used for the existing link in A2 (underline + blue font)
# by copying A2's formatting onto B2, then restore B2's own text/value.
$ws.Range("A2").Copy($ws.Range("B2"))
$ws.Range("B2").Value = $url
